# Updates the crypto price/volume table with the latest scraped values.
# Cells whose new text would otherwise be auto-parsed by Excel as a plain
# number (single dot, no surrounding spaces/'%') are forced to stay text
# via NumberFormat "@" before the value is assigned, so they keep exactly
# the same literal representation as the source data (e.g. "570.52" and
# not 570.52 as a float, "0.0000124" and not 1.24E-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.805.05'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '3.359.48'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.52'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.49'
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.357.68'
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E11").Value = '  -3.71%  '
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '3.934.18'
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.99'
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("E16").Value = '  -5.20%  '
$ws.Range("D17").Value = '3.357.45'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").Value = '60.986.32'
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.96'
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.80'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.25'
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.67'
$ws.Range("E22").Value = '  -4.38%  '
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("D24").Value = '3.500.81'
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.87'
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").Value = '  -3.96%  '
$ws.Range("E28").Value = '  +9.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.47'
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.166'
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.19'
$ws.Range("E36").Value = '  -6.14%  '
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("E38").Value = '  -3.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.70'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0749'
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.767'
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.69'
$ws.Range("E43").Value = '  -3.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.42'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.38'
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("E46").Value = '  -3.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.69'
$ws.Range("E47").Value = '  -6.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.14'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.78'
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").Value = '2.347.19'
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("E51").Value = '  +2.26%  '
